$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 4, duplicating the values from row 3 (Name/Image columns)
$ws.Range("A4").Value = $ws.Range("A3").Value2
$ws.Range("B4").Value = $ws.Range("B3").Value2

# Update selection to match the authored state
$ws.Range("B10").Select()
